# Sample Project / Main.xlsx - rule row 11 ("R40") now carries the literal
# text value "1" in its Rule-name cell (column B) instead of "R40".
#
# A leading apostrophe is the standard Excel way of forcing a number-looking
# entry ("1") to be stored as text rather than being auto-converted to a
# numeric value - exactly what happened here (the cell keeps its general
# number format, but the stored/shared-string value becomes the text "1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "'1"
